# Auto update Excel log
# Adds new sensor rows to mmWave, PIR, Humidity sheets,
# and creates two new sheets: Camera and Proximity.

$wb = $excel.ActiveWorkbook

function Set-RowText {
    param(
        $ws,
        [int]$row,
        [string[]]$values
    )
    # Force every cell in the row to be treated as plain text first, so
    # Excel's automatic date/time/percentage detection does not silently
    # convert values like "2026-01-30" or "86.9%" into numbers.
    $lastCol = $values.Length
    $ws.Range($ws.Cells.Item($row,1), $ws.Cells.Item($row,$lastCol)).NumberFormat = "@"
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# mmWave sheet: two new presence rows
# ---------------------------------------------------------------------------
$mmWave = $wb.Worksheets.Item("mmWave")
Set-RowText $mmWave 3 @("2026-01-30", "15:54:54", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
Set-RowText $mmWave 4 @("2026-01-30", "15:55:05", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")

# ---------------------------------------------------------------------------
# PIR sheet: six new "no motion" rows
# ---------------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")
Set-RowText $pir 3 @("2026-01-30", "15:54:52", "15:00", "Bathroom", "No Motion", "Inactive")
Set-RowText $pir 4 @("2026-01-30", "15:54:54", "15:00", "Bathroom", "No Motion", "Inactive")
Set-RowText $pir 5 @("2026-01-30", "15:54:57", "15:00", "Bathroom", "No Motion", "Inactive")
Set-RowText $pir 6 @("2026-01-30", "15:55:02", "15:00", "Bathroom", "No Motion", "Inactive")
Set-RowText $pir 7 @("2026-01-30", "15:55:07", "15:00", "Bathroom", "No Motion", "Inactive")
Set-RowText $pir 8 @("2026-01-30", "15:55:12", "15:00", "Bathroom", "No Motion", "Inactive")

# ---------------------------------------------------------------------------
# Humidity sheet: five new humidity reading rows
# ---------------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")
Set-RowText $humidity 3 @("2026-01-30", "15:54:52", "15:00", "Bathroom", "86.9%", "Active")
Set-RowText $humidity 4 @("2026-01-30", "15:54:54", "15:00", "Bathroom", "86.9%", "Active")
Set-RowText $humidity 5 @("2026-01-30", "15:54:57", "15:00", "Bathroom", "87.8%", "Active")
Set-RowText $humidity 6 @("2026-01-30", "15:55:02", "15:00", "Bathroom", "87.8%", "Active")
Set-RowText $humidity 7 @("2026-01-30", "15:55:07", "15:00", "Bathroom", "86.4%", "Active")

# ---------------------------------------------------------------------------
# New sheet: Camera (inserted after Humidity, the last existing sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$camera = $wb.Worksheets.Add($null, $lastSheet)
$camera.Name = "Camera"
Set-RowText $camera 1 @("Date", "Timestamp", "Hour", "Location", "Value", "Status")
Set-RowText $camera 2 @("2026-01-30", "15:54:55", "15:00", "Living Room Main Door", "Image Captured (EXIT)", "Active")
Set-RowText $camera 3 @("2026-01-30", "15:54:58", "15:00", "Living Room Main Door", "Image Captured (ENTER)", "Active")

# ---------------------------------------------------------------------------
# New sheet: Proximity (inserted after Camera)
# ---------------------------------------------------------------------------
$proximity = $wb.Worksheets.Add($null, $camera)
$proximity.Name = "Proximity"
Set-RowText $proximity 1 @("Date", "Timestamp", "Hour", "Location", "Value", "Status")
Set-RowText $proximity 2 @("2026-01-30", "15:54:55", "15:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
Set-RowText $proximity 3 @("2026-01-30", "15:54:58", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
